$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A4").Value = "102_AutomobileInsurance_003_InsurantData_001_MandatoryFields"
$ws.Range("C4").Value = "Insurant Page check for open mandatory fields"
$ws.Range("B4").Value = "Button Next from Page VehicleData"

$ws.Range("C7").Select()
